# Update Global_M2 "New Zealand_FX" sheet: correct the last existing monthly
# row and append three new monthly rows of FX_IDC:USDNZD data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the previously-last row (row 210): high/close values revised ---
$ws.Range("D210").Value = 1.6356
$ws.Range("F210").Value = 1.6172

# --- Append rows 211-213, copying formatting (date style, borders, etc.) ---
# from the existing last data row before overwriting the values.
$ws.Range("A210:G210").Copy($ws.Range("A211:G211"))
$ws.Range("A210:G210").Copy($ws.Range("A212:G212"))
$ws.Range("A210:G210").Copy($ws.Range("A213:G213"))

# Row 211 - 2023-05-01
$ws.Range("A211").Value = 45047.33333333334
$ws.Range("B211").Value = "FX_IDC:USDNZD"
$ws.Range("C211").Value = 1.6163
$ws.Range("D211").Value = 1.6702
$ws.Range("E211").Value = 1.5669
$ws.Range("F211").Value = 1.6614
$ws.Range("G211").Value = 0

# Row 212 - 2023-06-01
$ws.Range("A212").Value = 45078.33333333334
$ws.Range("B212").Value = "FX_IDC:USDNZD"
$ws.Range("C212").Value = 1.6614
$ws.Range("D212").Value = 1.6685
$ws.Range("E212").Value = 1.6013
$ws.Range("F212").Value = 1.6288
$ws.Range("G212").Value = 0

# Row 213 - 2023-07-03
$ws.Range("A213").Value = 45110.33333333334
$ws.Range("B213").Value = "FX_IDC:USDNZD"
$ws.Range("C213").Value = 1.6288
$ws.Range("D213").Value = 1.634
$ws.Range("E213").Value = 1.6082
$ws.Range("F213").Value = 1.6101
$ws.Range("G213").Value = 0
